# Auto-generated market data refresh for Sargatanas_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 10049.818
$ws.Range("J112").Value = 10380.762
$ws.Range("L112").Value = 31142.286
$ws.Range("N112").Value = -33358.286
$ws.Range("H132").Value = 1517.2285
$ws.Range("I132").Value = 1444.2059
$ws.Range("K132").Value = 4332.6177
$ws.Range("M132").Value = -1802.6177
$ws.Range("H135").Value = 238950.53
$ws.Range("I135").Value = 286543.5
$ws.Range("K135").Value = 2578891.5
$ws.Range("M135").Value = -2576356.5
$ws.Range("H138").Value = 3229617.2
$ws.Range("I138").Value = 2724.2856
$ws.Range("J138").Value = 5887058.5
$ws.Range("K138").Value = 8172.8568
$ws.Range("L138").Value = 17661175.5
$ws.Range("M138").Value = -3032.8568
$ws.Range("N138").Value = -17671455.5
$ws.Range("H141").Value = 3588.6667
$ws.Range("I141").Value = 3369.4546
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 10108.3638
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = -4928.363799999999
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4177.773
$ws.Range("I45").Value = 1862.8182
$ws.Range("K45").Value = 1862.8182
$ws.Range("M45").Value = -1485.8182
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H97").Value = 3339880.8
$ws.Range("I97").Value = 340.70587
$ws.Range("J97").Value = 10436403
$ws.Range("K97").Value = 340.70587
$ws.Range("L97").Value = 10436403
$ws.Range("M97").Value = 155.29413
$ws.Range("N97").Value = -10437395
$ws.Range("H122").Value = 2396.9033
$ws.Range("I122").Value = 1807.3684
$ws.Range("K122").Value = 5422.1052
$ws.Range("M122").Value = -2972.1052
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 59386
$ws.Range("J53").Value = 59386
$ws.Range("L53").Value = 59386
$ws.Range("N53").Value = -60534
$ws.Range("H105").Value = 3566.6191
$ws.Range("I105").Value = 2426.2727
$ws.Range("K105").Value = 2426.2727
$ws.Range("M105").Value = -679.2727
$ws.Range("H134").Value = 8599.538
$ws.Range("I134").Value = 2599.375
$ws.Range("K134").Value = 7798.125
$ws.Range("M134").Value = -5263.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3978.2646
$ws.Range("I16").Value = 3519.7273
$ws.Range("K16").Value = 3519.7273
$ws.Range("M16").Value = -3232.7273
$ws.Range("H18").Value = 39576.375
$ws.Range("J18").Value = 39576.375
$ws.Range("L18").Value = 39576.375
$ws.Range("N18").Value = -40036.375
$ws.Range("H53").Value = 55532.5
$ws.Range("J53").Value = 55532.5
$ws.Range("L53").Value = 55532.5
$ws.Range("N53").Value = -56746.5
$ws.Range("H58").Value = 7742.2607
$ws.Range("I58").Value = 1491.25
$ws.Range("J58").Value = 11076.134
$ws.Range("K58").Value = 1491.25
$ws.Range("L58").Value = 11076.134
$ws.Range("M58").Value = -1288.25
$ws.Range("N58").Value = -11482.134
$ws.Range("H113").Value = 3978.2646
$ws.Range("I113").Value = 3519.7273
$ws.Range("K113").Value = 3519.7273
$ws.Range("M113").Value = -1349.7273
$ws.Range("H136").Value = 7742.2607
$ws.Range("I136").Value = 1491.25
$ws.Range("J136").Value = 11076.134
$ws.Range("K136").Value = 4473.75
$ws.Range("L136").Value = 33228.402
$ws.Range("M136").Value = -1923.75
$ws.Range("N136").Value = -38328.402

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1824.6666
$ws.Range("I5").Value = 1778.0952
$ws.Range("J5").Value = 1889.8667
$ws.Range("K5").Value = 5334.2856
$ws.Range("L5").Value = 5669.6001
$ws.Range("M5").Value = -5222.2856
$ws.Range("N5").Value = -5893.6001
$ws.Range("H75").Value = 222223710
$ws.Range("I75").Value = 333333340
$ws.Range("J75").Value = 166668900
$ws.Range("K75").Value = 1000000020
$ws.Range("L75").Value = 500006700
$ws.Range("M75").Value = -999999022
$ws.Range("N75").Value = -500008696
$ws.Range("H78").Value = 222223710
$ws.Range("I78").Value = 333333340
$ws.Range("J78").Value = 166668900
$ws.Range("K78").Value = 3000000060
$ws.Range("L78").Value = 1500020100
$ws.Range("M78").Value = -2999995068
$ws.Range("N78").Value = -1500030084
$ws.Range("H113").Value = 1710
$ws.Range("J113").Value = 2164.75
$ws.Range("L113").Value = 6494.25
$ws.Range("N113").Value = -10834.25
$ws.Range("H117").Value = 150500340
$ws.Range("I117").Value = 27778074
$ws.Range("J117").Value = 334583700
$ws.Range("K117").Value = 83334222
$ws.Range("L117").Value = 1003751100
$ws.Range("M117").Value = -83330780
$ws.Range("N117").Value = -1003757984
$ws.Range("H131").Value = 2159
$ws.Range("I131").Value = 1831.5
$ws.Range("J131").Value = 2255.3235
$ws.Range("K131").Value = 5494.5
$ws.Range("L131").Value = 6765.970499999999
$ws.Range("M131").Value = -454.5
$ws.Range("N131").Value = -16845.9705
$ws.Range("H135").Value = 1824.6666
$ws.Range("I135").Value = 1778.0952
$ws.Range("J135").Value = 1889.8667
$ws.Range("K135").Value = 16002.8568
$ws.Range("L135").Value = 17008.8003
$ws.Range("M135").Value = -13467.8568
$ws.Range("N135").Value = -22078.8003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5317.769
$ws.Range("I126").Value = 2508.5
$ws.Range("K126").Value = 7525.5
$ws.Range("M126").Value = -5055.5
$ws.Range("H132").Value = 4620.9565
$ws.Range("I132").Value = 2691.4167
$ws.Range("J132").Value = 6725.909
$ws.Range("K132").Value = 8074.250100000001
$ws.Range("L132").Value = 20177.727
$ws.Range("M132").Value = -5544.250100000001
$ws.Range("N132").Value = -25237.727
$ws.Range("H134").Value = 87939.22
$ws.Range("J134").Value = 87939.22
$ws.Range("L134").Value = 263817.66
$ws.Range("N134").Value = -268887.66
$ws.Range("H141").Value = 59974.5
$ws.Range("J141").Value = 59974.5
$ws.Range("L141").Value = 59974.5
$ws.Range("N141").Value = -70334.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3298.7368
$ws.Range("I16").Value = 3216.2942
$ws.Range("K16").Value = 3216.2942
$ws.Range("M16").Value = -3046.2942
$ws.Range("H46").Value = 3099.4
$ws.Range("I46").Value = 1373
$ws.Range("J46").Value = 4250.3335
$ws.Range("K46").Value = 1373
$ws.Range("L46").Value = 4250.3335
$ws.Range("M46").Value = -1185
$ws.Range("N46").Value = -4626.3335
$ws.Range("H55").Value = 944.8
$ws.Range("I55").Value = 1010.2222
$ws.Range("K55").Value = 1010.2222
$ws.Range("M55").Value = -837.2222
$ws.Range("H56").Value = 500
$ws.Range("I56").Value = 500
$ws.Range("K56").Value = 500
$ws.Range("H122").Value = 8804.723
$ws.Range("I122").Value = 11069.857
$ws.Range("J122").Value = 7363.273
$ws.Range("K122").Value = 33209.571
$ws.Range("L122").Value = 22089.819
$ws.Range("M122").Value = -30759.571
$ws.Range("N122").Value = -26989.819
$ws.Range("H132").Value = 6855.5293
$ws.Range("I132").Value = 4050.2144
$ws.Range("K132").Value = 12150.6432
$ws.Range("M132").Value = -9620.643199999999
$ws.Range("H139").Value = 78568
$ws.Range("J139").Value = 78568
$ws.Range("L139").Value = 78568
$ws.Range("N139").Value = -88848
$ws.Range("M56").Value = 191

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 81666.586
$ws.Range("I29").Value = 88181.82000000001
$ws.Range("K29").Value = 88181.82000000001
$ws.Range("M29").Value = -87891.82000000001
$ws.Range("H122").Value = 9339637
$ws.Range("I122").Value = 13625790
$ws.Range("J122").Value = 10950.765
$ws.Range("K122").Value = 40877370
$ws.Range("L122").Value = 32852.295
$ws.Range("M122").Value = -40874920
$ws.Range("N122").Value = -37752.295
$ws.Range("H126").Value = 4118.1
$ws.Range("I126").Value = 1896.4546
$ws.Range("K126").Value = 5689.3638
$ws.Range("M126").Value = -3219.3638
$ws.Range("H133").Value = 219983.33
$ws.Range("J133").Value = 219983.33
$ws.Range("L133").Value = 219983.33
$ws.Range("N133").Value = -230103.33
$ws.Range("H140").Value = 74380
$ws.Range("J140").Value = 74380
$ws.Range("L140").Value = 74380
$ws.Range("N140").Value = -84740
$ws.Range("H141").Value = 66649.336
$ws.Range("J141").Value = 84974
$ws.Range("L141").Value = 84974
$ws.Range("N141").Value = -95334
